{"js": "// Replace the multiplication-problem text in each table cell with the\n// updated operands, per the commit's regenerated number set.\n// Each \"old\" string occurs exactly once in the document, so a scoped\n// search + full-match replace is safe and keeps existing run formatting\n// (font, size, etc.) untouched.\nconst replacements = [\n  [\"26\u00d738=\", \"29\u00d769=\"],\n  [\"97\u00d793=\", \"19\u00d723=\"],\n  [\"66\u00d796=\", \"64\u00d749=\"],\n  [\"98\u00d747=\", \"68\u00d713=\"],\n  [\"43\u00d726=\", \"30\u00d754=\"],\n  [\"67\u00d768=\", \"93\u00d783=\"],\n  [\"73\u00d734=\", \"93\u00d783=\"],\n  [\"53\u00d740=\", \"86\u00d764=\"],\n  [\"71\u00d721=\", \"49\u00d712=\"],\n  [\"43\u00d797=\", \"83\u00d778=\"],\n  [\"90\u00d773=\", \"87\u00d738=\"],\n  [\"11\u00d734=\", \"58\u00d743=\"],\n  [\"70\u00d782=\", \"92\u00d761=\"],\n  [\"26\u00d714=\", \"26\u00d791=\"],\n  [\"85\u00d798=\", \"65\u00d773=\"],\n  [\"71\u00d711=\", \"59\u00d735=\"],\n  [\"42\u00d777=\", \"31\u00d728=\"],\n  [\"56\u00d770=\", \"72\u00d756=\"],\n  [\"47\u00d743=\", \"93\u00d767=\"],\n  [\"55\u00d743=\", \"42\u00d743=\"],\n  [\"94\u00d773=\", \"35\u00d716=\"],\n  [\"42\u00d740=\", \"68\u00d784=\"],\n  [\"38\u00d794=\", \"83\u00d728=\"],\n  [\"25\u00d721=\", \"90\u00d796=\"],\n  [\"36\u00d768=\", \"68\u00d735=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the multiplication-problem text in each table cell with the\n# regenerated operand pairs. Each \"old\" string occurs exactly once in\n# the document, so Find/Replace scoped to the whole document content is\n# safe and leaves the surrounding run formatting (font, size) intact.\n# Parallel arrays (rather than an array-of-arrays) keep the index math\n# simple and avoid any ambiguity around nested-array literals.\n\n$d = $word.ActiveDocument\n\n$oldTexts = @(\"26\u00d738=\", \"97\u00d793=\", \"66\u00d796=\", \"98\u00d747=\", \"43\u00d726=\", \"67\u00d768=\", \"73\u00d734=\", \"53\u00d740=\", \"71\u00d721=\", \"43\u00d797=\", \"90\u00d773=\", \"11\u00d734=\", \"70\u00d782=\", \"26\u00d714=\", \"85\u00d798=\", \"71\u00d711=\", \"42\u00d777=\", \"56\u00d770=\", \"47\u00d743=\", \"55\u00d743=\", \"94\u00d773=\", \"42\u00d740=\", \"38\u00d794=\", \"25\u00d721=\", \"36\u00d768=\")\n$newTexts = @(\"29\u00d769=\", \"19\u00d723=\", \"64\u00d749=\", \"68\u00d713=\", \"30\u00d754=\", \"93\u00d783=\", \"93\u00d783=\", \"86\u00d764=\", \"49\u00d712=\", \"83\u00d778=\", \"87\u00d738=\", \"58\u00d743=\", \"92\u00d761=\", \"26\u00d791=\", \"65\u00d773=\", \"59\u00d735=\", \"31\u00d728=\", \"72\u00d756=\", \"93\u00d767=\", \"42\u00d743=\", \"35\u00d716=\", \"68\u00d784=\", \"83\u00d728=\", \"90\u00d796=\", \"68\u00d735=\")\n\nfor ($i = 0; $i -lt $oldTexts.Count; $i++) {\n    $oldText = $oldTexts[$i]\n    $newText = $newTexts[$i]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
